$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.918.34"
$ws.Range("E2").Value = "  +3.68%  "
$ws.Range("D3").Value = "'2.592.74"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'521.80"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "'140.99"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "'2.616.20"
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").Value = "'6.52"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("D14").Value = "'3.056.84"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "'58.934.43"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("D16").Value = "'20.55"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.620.49"
$ws.Range("E17").Value = "  +3.03%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'339.92"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").Value = "'4.33"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").Value = "'10.22"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").Value = "'6.52"
$ws.Range("E22").Value = "  +6.20%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "'66.23"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'7.14"
$ws.Range("E28").Value = "  +4.06%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "'148.93"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'4.01"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D37").Value = "'36.34"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").Value = "'0.838"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "'3.55"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "'277.52"
$ws.Range("E42").Value = "  +6.35%  "
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "'0.591"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("D48").Value = "'18.60"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "'1.986.43"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0221"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.51"
$ws.Range("E51").Value = "  +0.30%  "
